$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "航天发展"
$ws.Range("C2").Value = "航天发展"
$ws.Range("A3").Value = "平潭发展"
$ws.Range("C3").Value = "平潭发展"
$ws.Range("A4").Value = "航天动力"
$ws.Range("B4").Value = "N中国铀"
$ws.Range("C4").Value = "中国铀业"
$ws.Range("A5").Value = "N中国铀"
$ws.Range("B5").Value = "航天动力"
$ws.Range("C5").Value = "海王生物"
$ws.Range("A6").Value = "道明光学"
$ws.Range("B6").Value = "合富中国"
$ws.Range("C6").Value = "实达集团"
$ws.Range("A7").Value = "实达集团"
$ws.Range("B7").Value = "海欣食品"
$ws.Range("B8").Value = "海王生物"
$ws.Range("C8").Value = "金富科技"
$ws.Range("A9").Value = "福蓉科技"
$ws.Range("B9").Value = "雷科防务"
$ws.Range("C9").Value = "道明光学"
$ws.Range("A10").Value = "海王生物"
$ws.Range("B10").Value = "实达集团"
$ws.Range("C10").Value = "航天动力"
$ws.Range("A11").Value = "海欣食品"
$ws.Range("B11").Value = "福蓉科技"
$ws.Range("C11").Value = "华映科技"
$ws.Range("A12").Value = "赛微电子"
$ws.Range("B12").Value = "道明光学"
$ws.Range("C12").Value = "福蓉科技"
$ws.Range("A13").Value = "雷科防务"
$ws.Range("B13").Value = "华映科技"
$ws.Range("C13").Value = "合富中国"
$ws.Range("A14").Value = "顺灏股份"
$ws.Range("B14").Value = "赛微电子"
$ws.Range("C14").Value = "雷科防务"
$ws.Range("A15").Value = "华映科技"
$ws.Range("B15").Value = "安泰集团"
$ws.Range("C15").Value = "安泰集团"
$ws.Range("A16").Value = "安泰集团"
$ws.Range("B16").Value = "海马汽车"
$ws.Range("C16").Value = "通宇通讯"
$ws.Range("A17").Value = "金富科技"
$ws.Range("C17").Value = "银河电子"
$ws.Range("A18").Value = "海马汽车"
$ws.Range("B18").Value = "顺灏股份"
$ws.Range("C18").Value = "赛微电子"
$ws.Range("A19").Value = "榕基软件"
$ws.Range("B19").Value = "黄河旋风"
$ws.Range("C19").Value = "顺灏股份"
$ws.Range("A20").Value = "通宇通讯"
$ws.Range("B20").Value = "大有能源"
$ws.Range("C20").Value = "安记食品"
$ws.Range("A21").Value = "襄阳轴承"
$ws.Range("B21").Value = "银河电子"
$ws.Range("C21").Value = "榕基软件"
